$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Fitness) values for rows 2-124 according to the new run's data.
$ws.Range("C2:C26").Value = 7320
$ws.Range("C27:C55").Value = 7318
$ws.Range("C56:C97").Value = 7310
$ws.Range("C98:C124").Value = 7295
